$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in column D was renamed from "gender" to "sex".
$ws.Range("D1").Value = "sex"
$ws.Range("D1").Select() | Out-Null
